# Add a new "Website - Menulayout" entry to the "03_Grießer" sheet and make
# that sheet the active tab (matching the author's last-saved UI state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("03_Grießer")

# Copy the formatting of the last existing row down into the new row so the
# new cells pick up the same date / percentage number formats as the rest of
# the table, then activate the sheet (this is what drives tabSelected /
# activeTab on save) and fill in the new values.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

$ws.Activate()

$ws.Range("A8").Value = 42660
$ws.Range("B8").Value = "Website - Menulayout"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 30

$ws.Range("B8").Select()
